{"js": "// Find the paragraph containing \".background(Color.Red)\" (the one under the\n// \"Modifier Attributes\" heading) and insert two new paragraphs right after it:\n//   1. \".fillMaxHeight(), .fillMaxWidth(), .width(), .height(), .fillMaxSize()\"\n//   2. \"You can pass paremeters to Maximum height, width and fill maxsize like 0.5f 50% covering.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \".background(Color.Red)\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find paragraph '.background(Color.Red)'\");\n}\n\n// Insert the second new line first (directly \"After\" pushes subsequent\n// inserts further down), then insert the first line right after the target\n// paragraph so the final order matches the diff.\nconst secondPara = target.insertParagraph(\n  \"You can pass paremeters to Maximum height, width and fill maxsize like 0.5f 50% covering.\",\n  \"After\"\n);\nconst firstPara = target.insertParagraph(\n  \".fillMaxHeight(), .fillMaxWidth(), .width(), .height(), .fillMaxSize()\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Find the paragraph containing \".background(Color.Red)\" (the one under the\n# \"Modifier Attributes\" heading) and insert two new paragraphs right after it:\n#   1. \".fillMaxHeight(), .fillMaxWidth(), .width(), .height(), .fillMaxSize()\"\n#   2. \"You can pass paremeters to Maximum height, width and fill maxsize like 0.5f 50% covering.\"\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \".background(Color.Red)\") {\n        $target = $p\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph '.background(Color.Red)'\"\n}\n\n# Insert the first new paragraph right after the target paragraph, then the\n# second new paragraph right after the first, so the final order matches the\n# diff.\n$target.Range.InsertParagraphAfter()\n$firstNew = $target.Next()\n$firstNew.Range.Text = \".fillMaxHeight(), .fillMaxWidth(), .width(), .height(), .fillMaxSize()\"\n\n$firstNew.Range.InsertParagraphAfter()\n$secondNew = $firstNew.Next()\n$secondNew.Range.Text = \"You can pass paremeters to Maximum height, width and fill maxsize like 0.5f 50% covering.\"\n"}
